$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8502226508649073
$ws.Range("C2").Value = 0.274295890186977
$ws.Range("D2").Value = 0.02262213656532452
$ws.Range("F2").Value = 0.6855360736699083
$ws.Range("G2").Value = 0.00241647875992391
$ws.Range("I2").Value = 0.6316091717484298
$ws.Range("L2").Value = 0.2797379928387187
$ws.Range("M2").Value = 0.2204871106269835
$ws.Range("N2").Value = 1.263930770199558
$ws.Range("O2").Value = 2.321387632367447

$ws.Range("B3").Value = 0.7678603602510066
$ws.Range("C3").Value = 0.2612004100952277
$ws.Range("D3").Value = 0.02115128993445836
$ws.Range("F3").Value = 0.6798511967606018
$ws.Range("G3").Value = 0.002419142853369471
$ws.Range("I3").Value = 0.6368723246262959
$ws.Range("L3").Value = 0.277824218253258
$ws.Range("M3").Value = 0.2073207977430371
$ws.Range("N3").Value = 1.274956989684959
$ws.Range("O3").Value = 2.316338765796019

$ws.Range("B4").Value = 0.717367556390883
$ws.Range("C4").Value = 0.2530901214990422
$ws.Range("D4").Value = 0.02024077513966205
$ws.Range("F4").Value = 0.6768320865624275
$ws.Range("G4").Value = 0.002420867562895975
$ws.Range("I4").Value = 0.6404838703671061
$ws.Range("L4").Value = 0.2767955168151985
$ws.Range("M4").Value = 0.1993093642033941
$ws.Range("N4").Value = 1.282208482665915
$ws.Range("O4").Value = 2.314797927902845

$ws.Range("B5").Value = 0.6968122334308191
$ws.Range("C5").Value = 0.2497677812238663
$ws.Range("D5").Value = 0.01986789001968958
$ws.Range("F5").Value = 0.6757202660513855
$ws.Range("G5").Value = 0.002421592828447503
$ws.Range("I5").Value = 0.642051084127722
$ws.Range("L5").Value = 0.2764131683536846
$ws.Range("M5").Value = 0.1960631333633529
$ws.Range("N5").Value = 1.285284647921792
$ws.Range("O5").Value = 2.314561833932203

$ws.Range("B6").Value = 0.6934003385159144
$ws.Range("C6").Value = 0.2492150674320044
$ws.Range("D6").Value = 0.01980586201222678
$ws.Range("F6").Value = 0.6755428041529257
$ws.Range("G6").Value = 0.002421714614976039
$ws.Range("I6").Value = 0.6423170837063878
$ws.Range("L6").Value = 0.2763519071917599
$ws.Range("M6").Value = 0.1955252221524404
$ws.Range("N6").Value = 1.285802760880927
$ws.Range("O6").Value = 2.314546287472183

$ws.Range("B7").Value = 0.7170902537554582
$ws.Range("C7").Value = 0.2530453852019008
$ws.Range("D7").Value = 0.02023575371608644
$ws.Range("F7").Value = 0.6768166124659913
$ws.Range("G7").Value = 0.002420877253194516
$ws.Range("I7").Value = 0.6405046198493167
$ws.Range("L7").Value = 0.2767902110340685
$ws.Range("M7").Value = 0.1992655092503099
$ws.Range("N7").Value = 1.282249478371497
$ws.Range("O7").Value = 2.31479315778202

$ws.Range("B8").Value = 0.8218087011202897
$ws.Range("C8").Value = 0.2697951161098047
$ws.Range("D8").Value = 0.02211653804376823
$ws.Range("F8").Value = 0.6834780478668065
$ws.Range("G8").Value = 0.002417378918612267
$ws.Range("I8").Value = 0.6333450112017971
$ws.Range("L8").Value = 0.2790477750013949
$ws.Range("M8").Value = 0.2159323871763519
$ws.Range("N8").Value = 1.267632756332929
$ws.Range("O8").Value = 2.319323008878769

$ws.Range("B9").Value = 1.027733164201663
$ws.Range("C9").Value = 0.3020828056101834
$ws.Range("D9").Value = 0.02574523859521349
$ws.Range("F9").Value = 0.7002858259977671
$ws.Range("G9").Value = 0.002411221411438744
$ws.Range("I9").Value = 0.6223227679381829
$ws.Range("L9").Value = 0.2846345087666435
$ws.Range("M9").Value = 0.2491860462005704
$ws.Range("N9").Value = 1.242785025059348
$ws.Range("O9").Value = 2.340593385131939

$ws.Range("B10").Value = 1.179326636567737
$ws.Range("C10").Value = 0.3254579366024188
$ws.Range("D10").Value = 0.02837422519262844
$ws.Range("F10").Value = 0.7149257070329185
$ws.Range("G10").Value = 0.002407121650231514
$ws.Range("I10").Value = 0.6160689020954919
$ws.Range("L10").Value = 0.2894446495395613
$ws.Range("M10").Value = 0.2739579567332413
$ws.Range("N10").Value = 1.226850552802254
$ws.Range("O10").Value = 2.363800422514771

$ws.Range("B11").Value = 1.248346503789435
$ws.Range("C11").Value = 0.3360154372586237
$ws.Range("D11").Value = 0.02956203724794193
$ws.Range("F11").Value = 0.7220852153004671
$ws.Range("G11").Value = 0.002405347764396638
$ws.Range("I11").Value = 0.6136252960728115
$ws.Range("L11").Value = 0.2917858744322359
$ws.Range("M11").Value = 0.2852999361473394
$ws.Range("N11").Value = 1.220104599513853
$ws.Range("O11").Value = 2.376010148408767

$ws.Range("B12").Value = 1.274489932431095
$ws.Range("C12").Value = 0.3400022158864715
$ws.Range("D12").Value = 0.03001064444065094
$ws.Range("F12").Value = 0.724868300068195
$ws.Range("G12").Value = 0.002404689074276644
$ws.Range("I12").Value = 0.612757752504038
$ws.Range("L12").Value = 0.2926944083428111
$ws.Range("M12").Value = 0.2896051914669258
$ws.Range("N12").Value = 1.217622313694442
$ws.Range("O12").Value = 2.380871710214336

$ws.Range("B13").Value = 1.268859183921563
$ws.Range("C13").Value = 0.3391440893315973
$ws.Range("D13").Value = 0.02991408220839986
$ws.Range("F13").Value = 0.7242657125612482
$ws.Range("G13").Value = 0.002404830355934832
$ws.Range("I13").Value = 0.6129420215786325
$ws.Range("L13").Value = 0.2924977632735875
$ws.Range("M13").Value = 0.2886775233778422
$ws.Range("N13").Value = 1.218153704998606
$ws.Range("O13").Value = 2.37981409519017

$ws.Range("B14").Value = 1.250497205146985
$ws.Range("C14").Value = 0.3363436558651358
$ws.Range("D14").Value = 0.02959896844116372
$ws.Range("F14").Value = 0.7223127392892366
$ws.Range("G14").Value = 0.00240529331246933
$ws.Range("I14").Value = 0.6135527638930967
$ws.Range("L14").Value = 0.2918601801716818
$ws.Range("M14").Value = 0.2856539270473277
$ws.Range("N14").Value = 1.219898932319275
$ws.Range("O14").Value = 2.37640534012391

$ws.Range("B15").Value = 1.239250845982269
$ws.Range("C15").Value = 0.3346268558847783
$ws.Range("D15").Value = 0.02940579632086582
$ws.Range("F15").Value = 0.7211258568127192
$ws.Range("G15").Value = 0.002405578583386016
$ws.Range("I15").Value = 0.613934390715059
$ws.Range("L15").Value = 0.2914725008934624
$ws.Range("M15").Value = 0.2838032210364645
$ws.Range("N15").Value = 1.220977344046176
$ws.Range("O15").Value = 2.374348384443664

$ws.Range("B16").Value = 1.17481709689099
$ws.Range("C16").Value = 0.3247664339102698
$ws.Range("D16").Value = 0.02829643339497778
$ws.Range("F16").Value = 0.714467878699665
$ws.Range("G16").Value = 0.002407239405997988
$ws.Range("I16").Value = 0.6162366804946373
$ws.Range("L16").Value = 0.2892947209931691
$ws.Range("M16").Value = 0.2732181841441488
$ws.Range("N16").Value = 1.22730153003539
$ws.Range("O16").Value = 2.363035774752262

$ws.Range("B17").Value = 1.135303243149565
$ws.Range("C17").Value = 0.3186977899926262
$ws.Range("D17").Value = 0.02761377741199311
$ws.Range("F17").Value = 0.7105114741620895
$ws.Range("G17").Value = 0.002408281560882965
$ws.Range("I17").Value = 0.6177519091526449
$ws.Range("L17").Value = 0.2879978937086918
$ws.Range("M17").Value = 0.2667431765268802
$ws.Range("N17").Value = 1.231309948877097
$ws.Range("O17").Value = 2.356519403797023

$ws.Range("B18").Value = 1.112581579667051
$ws.Range("C18").Value = 0.3152001280583931
$ws.Range("D18").Value = 0.02722036814400042
$ws.Range("F18").Value = 0.7082828899102651
$ws.Range("G18").Value = 0.002408889560581717
$ws.Range("I18").Value = 0.6186611963738819
$ws.Range("L18").Value = 0.2872664033008192
$ws.Range("M18").Value = 0.2630258182791678
$ws.Range("N18").Value = 1.233662805111535
$ws.Range("O18").Value = 2.352926901852385

$ws.Range("B19").Value = 1.104889430101764
$ws.Range("C19").Value = 0.3140146597108071
$ws.Range("D19").Value = 0.02708703606988649
$ws.Range("F19").Value = 0.707536405708268
$ws.Range("G19").Value = 0.002409096894244638
$ws.Range("I19").Value = 0.6189755493316973
$ws.Range("L19").Value = 0.2870212096709679
$ws.Range("M19").Value = 0.2617683755468505
$ws.Range("N19").Value = 1.234467569995225
$ws.Range("O19").Value = 2.351737246448295

$ws.Range("B20").Value = 1.139508984794986
$ws.Range("C20").Value = 0.3193445477097043
$ws.Range("D20").Value = 0.02768652650152603
$ws.Range("F20").Value = 0.7109277717176212
$ws.Range("G20").Value = 0.002408169734219145
$ws.Range("I20").Value = 0.6175867010040079
$ws.Range("L20").Value = 0.2881344521181859
$ws.Range("M20").Value = 0.2674317401393793
$ws.Range("N20").Value = 1.230878348661683
$ws.Range("O20").Value = 2.357196982066966

$ws.Range("B21").Value = 1.255890381915606
$ws.Range("C21").Value = 0.3371665143361611
$ws.Range("D21").Value = 0.02969155750970032
$ws.Range("F21").Value = 0.7228844220477839
$ws.Range("G21").Value = 0.002405156977515522
$ws.Range("I21").Value = 0.6133718047668921
$ws.Range("L21").Value = 0.2920468581023528
$ws.Range("M21").Value = 0.2865417529606091
$ws.Range("N21").Value = 1.219384355846202
$ws.Range("O21").Value = 2.377400112048008

$ws.Range("B22").Value = 1.331993097972372
$ws.Range("C22").Value = 0.3487493070726089
$ws.Range("D22").Value = 0.03099500589316051
$ws.Range("F22").Value = 0.731118104382432
$ws.Range("G22").Value = 0.002403263957614366
$ws.Range("I22").Value = 0.6109540527754689
$ws.Range("L22").Value = 0.2947318354725468
$ws.Range("M22").Value = 0.2990911496828375
$ws.Range("N22").Value = 1.212293544220714
$ws.Range("O22").Value = 2.391991449397551

$ws.Range("B23").Value = 1.291372385149486
$ws.Range("C23").Value = 0.3425733544109164
$ws.Range("D23").Value = 0.03029997490061476
$ws.Range("F23").Value = 0.7266852424540957
$ws.Range("G23").Value = 0.002404267364756139
$ws.Range("I23").Value = 0.6122135943156621
$ws.Range("L23").Value = 0.2932871170031035
$ws.Range("M23").Value = 0.2923878923189918
$ws.Range("N23").Value = 1.216039514978121
$ws.Range("O23").Value = 2.38407671611094

$ws.Range("B24").Value = 1.137607583287831
$ws.Range("C24").Value = 0.3190521756864655
$ws.Range("D24").Value = 0.02765363956785194
$ws.Range("F24").Value = 0.7107394202943027
$ws.Range("G24").Value = 0.002408220263552764
$ws.Range("I24").Value = 0.6176612727414437
$ws.Range("L24").Value = 0.2880726702229879
$ws.Range("M24").Value = 0.267120424285487
$ws.Range("N24").Value = 1.231073324467516
$ws.Range("O24").Value = 2.356890169732566

$ws.Range("B25").Value = 0.9719688714693007
$ws.Range("C25").Value = 0.2934085580656642
$ws.Range("D25").Value = 0.02477002653551352
$ws.Range("F25").Value = 0.6953371131148387
$ws.Range("G25").Value = 0.002412812389209757
$ws.Range("I25").Value = 0.6249810438009114
$ws.Range("L25").Value = 0.2829991433425079
$ws.Range("M25").Value = 0.2401297579099548
$ws.Range("N25").Value = 1.249099034058744
$ws.Range("O25").Value = 2.333510126506582

